# Adds the forecasted quarter/week columns' figures by updating the cached
# projection values on the "Supply_Demand" and "Wafer Plan" sheets.
$wb = $excel.ActiveWorkbook

# --- Supply_Demand: Yielded Supply / Inventory Balance rows for product 22B/23C (rows 8,12,13,14,18,19) ---
$wsSupplyDemand = $wb.Worksheets.Item("Supply_Demand")

$supplyDemandUpdates = @(
    ,("E8", 3273593400)
    ,("F8", 3273759304.805186)
    ,("G8", 3273837162.720054)
    ,("J8", 2046149837.078962)
    ,("K8", 2046156718.431006)
    ,("L8", 2046159693.931045)
    ,("N8", 2046169125)
    ,("E12", 2893928385.25)
    ,("F12", 3854856015.755187)
    ,("G12", 5886829119.099483)
    ,("H12", 5272886843.530636)
    ,("I12", 5887565329.13924)
    ,("J12", 6453908092.180819)
    ,("K12", 7167310482.178673)
    ,("L12", 8029655712.70414)
    ,("M12", 7119481700.387333)
    ,("N12", 8315355358.094958)
    ,("O12", 10143797197.98009)
    ,("P12", 12018178480.22488)
    ,("Q12", 13943645194.59819)
    ,("R12", 17045105823.81886)
    ,("S12", 20226140285.70691)
    ,("T12", 21470057933.70691)
    ,("U12", 22317906502.70691)
    ,("V12", 21971743057.70691)
    ,("W12", 22317739618.70691)
    ,("X12", 21805179552.70691)
    ,("Y12", 24304630383.70691)
    ,("Z12", 26343182464.70691)
    ,("AA12", 28565066345.70691)
    ,("AB12", 28330917656.70691)
    ,("AC12", 28660932822.70691)
    ,("AD12", 30891246292.70691)
    ,("AE12", 32020131634.70691)
    ,("AF12", 32174067894.70691)
    ,("AG12", 33407205964.70691)
    ,("AH12", 34717210734.70691)
    ,("AI12", 36359490002.70691)
    ,("AJ12", 37602017975.70691)
    ,("AK12", 38748807217.70691)
    ,("E13", -4049998001.956846)
    ,("F13", -2905690044.040411)
    ,("G13", -889000073.676734)
    ,("H13", 908016894.9958746)
    ,("I13", 1244767893.018836)
    ,("J13", 2132989790.914745)
    ,("K13", 2871927333.517089)
    ,("L13", 3768294393.348852)
    ,("M13", 2881293924.282703)
    ,("N13", 4185668876.667447)
    ,("O13", 6021997168.090668)
    ,("P13", 7905162302.698104)
    ,("Q13", 9830635347.961536)
    ,("E14", 3324968346.84)
    ,("F14", 2911621667.5)
    ,("G14", 2754235586.22475)
    ,("J14", 393461932.8486817)
    ,("K14", 393462387.4999999)
    ,("L14", 393462387.4999999)
    ,("N14", 393462387.4999999)
    ,("E18", -5625449407.817211)
    ,("F18", -8580904651.797619)
    ,("G18", -12103249703.41854)
    ,("H18", -23623922171.50393)
    ,("I18", -29111522546.16782)
    ,("J18", -33836086246.9054)
    ,("K18", -36231833151.4959)
    ,("L18", -37955283175.19156)
    ,("M18", -42502060339.13037)
    ,("N18", -43605526087.89382)
    ,("O18", -44688691675.31496)
    ,("P18", -45813667469.05714)
    ,("Q18", -46554027102.91869)
    ,("R18", -44699404399.0709)
    ,("S18", -42580469219.38747)
    ,("T18", -48908537796.88747)
    ,("U18", -50150407555.63747)
    ,("V18", -58401721569.13747)
    ,("W18", -71267329534.63747)
    ,("X18", -85510833105.13747)
    ,("Y18", -87447412682.48747)
    ,("Z18", -87612448123.98747)
    ,("AA18", -91195047282.48747)
    ,("AB18", -92529786118.23747)
    ,("AC18", -94188297366.78748)
    ,("AD18", -100162793044.2875)
    ,("AE18", -108369117547.7875)
    ,("AF18", -114387887964.5375)
    ,("AG18", -125096561788.0875)
    ,("AH18", -133473314510.5875)
    ,("AI18", -134455755482.0875)
    ,("AJ18", -135252765776.8375)
    ,("AK18", -138214717235.3875)
    ,("E19", -13223144679.50558)
    ,("F19", -15659464114.36675)
    ,("G19", -19117236579.2967)
    ,("H19", -31938038411.25397)
    ,("I19", -31538283345.12726)
    ,("J19", -35180851971.02086)
    ,("K19", -37442140408.73502)
    ,("L19", -39071557253.62081)
    ,("M19", -43588370741.38306)
    ,("N19", -44687776457.87805)
    ,("O19", -45779304086.56339)
    ,("P19", -46827356648.32944)
    ,("Q19", -47547050364.75099)
)

foreach ($update in $supplyDemandUpdates) {
    $wsSupplyDemand.Range($update[0]).Value = $update[1]
}

# --- Wafer Plan: weekly wafer-start plan rows for 22B (row 5) and 23C (row 6) ---
$wsWaferPlan = $wb.Worksheets.Item("Wafer Plan")

$waferPlanUpdates = @(
    ,("O5", 4000)
    ,("P5", 4000)
    ,("Q5", 4000)
    ,("R5", 4000)
    ,("S5", 4000)
    ,("T5", 4000)
    ,("U5", 4000)
    ,("V5", 4000)
    ,("W5", 4000)
    ,("X5", 4000)
    ,("Y5", 4000)
    ,("Z5", 4000)
    ,("AA5", 4000)
    ,("AB5", 4000)
    ,("AC5", 4000)
    ,("AD5", 4000)
    ,("AE5", 4000)
    ,("AF5", 4000)
    ,("AG5", 4000)
    ,("AH5", 4000)
    ,("AI5", 4000)
    ,("AO5", 4000)
    ,("AP5", 4000)
    ,("AQ5", 4000)
    ,("AR5", 4000)
    ,("AS5", 4000)
    ,("AT5", 4000)
    ,("AU5", 4000)
    ,("AV5", 4000)
    ,("AW5", 4000)
    ,("AX5", 4000)
    ,("AZ5", 4000)
    ,("BA5", 4000)
    ,("CD5", 2500)
    ,("CE5", 2500)
    ,("CF5", 2500)
    ,("CG5", 2500)
    ,("CH5", 2500)
    ,("CI5", 2500)
    ,("CJ5", 2500)
    ,("CP5", 2500)
    ,("CQ5", 2500)
    ,("CR5", 2500)
    ,("CS5", 2500)
    ,("CT5", 2500)
    ,("CU5", 2500)
    ,("CV5", 2500)
    ,("CW5", 2500)
    ,("DC5", 2500)
    ,("DD5", 2500)
    ,("DE5", 2500)
    ,("DF5", 2500)
    ,("DG5", 2500)
    ,("DH5", 2500)
    ,("DI5", 2500)
    ,("DJ5", 2500)
    ,("DL5", 2500)
    ,("EF5", 2500)
    ,("EG5", 2500)
    ,("EH5", 2500)
    ,("EI5", 2500)
    ,("EJ5", 2500)
    ,("O6", 2800)
    ,("P6", 3600)
    ,("Q6", 4400)
    ,("R6", 4570)
    ,("S6", 3718)
    ,("T6", 3169)
    ,("U6", 4501)
    ,("V6", 4700)
    ,("W6", 4700)
    ,("X6", 4700)
    ,("Y6", 4700)
    ,("Z6", 4700)
    ,("AA6", 4700)
    ,("AB6", 3700)
    ,("AC6", 3700)
    ,("AD6", 3700)
    ,("AE6", 3700)
    ,("AF6", 3700)
    ,("AG6", 3700)
    ,("AH6", 3700)
    ,("AI6", 3700)
    ,("AO6", 3500)
    ,("AP6", 3500)
    ,("AQ6", 3500)
    ,("AR6", 3500)
    ,("AS6", 3500)
    ,("AT6", 3500)
    ,("AU6", 3500)
    ,("AV6", 3500)
    ,("AW6", 3500)
    ,("AX6", 3500)
    ,("AZ6", 3500)
    ,("BA6", 3500)
    ,("CD6", 500)
    ,("CE6", 500)
    ,("CF6", 500)
    ,("CG6", 500)
    ,("CH6", 500)
    ,("CI6", 500)
    ,("CJ6", 500)
    ,("CP6", 500)
    ,("CQ6", 500)
    ,("CR6", 500)
    ,("CS6", 500)
    ,("CT6", 500)
    ,("CU6", 500)
    ,("CV6", 500)
    ,("CW6", 500)
    ,("DC6", 500)
    ,("DD6", 500)
    ,("DE6", 500)
    ,("DF6", 500)
    ,("DG6", 500)
    ,("DH6", 500)
    ,("DI6", 500)
    ,("DJ6", 500)
    ,("DL6", 500)
    ,("EF6", 500)
    ,("EG6", 500)
    ,("EH6", 500)
    ,("EI6", 500)
    ,("EJ6", 500)
)

foreach ($update in $waferPlanUpdates) {
    $wsWaferPlan.Range($update[0]).Value = $update[1]
}

Write-Output ("Updated " + $supplyDemandUpdates.Length + " Supply_Demand cells and " + $waferPlanUpdates.Length + " Wafer Plan cells")
